# Time series plots are staged.
# Adds BLS unemployment-series notes, a Portland population-history citation,
# and a BLS recession-article citation to the bibliography.

$d = $word.ActiveDocument

# Use explicit string literals (backtick escapes) for control characters so
# PowerShell never treats two adjacent [char] values as a numeric addition.
$cr = "`r"

# ---------------------------------------------------------------------
# Edit 1: expand the ggmap/BLS paragraph that currently reads
#   <br/>Datasets. (n.d.). Retrieved March 28, 2017, from ...
# into a BLS citation block followed by the original "Datasets." text in
# its own paragraph.
# ---------------------------------------------------------------------
$range = $d.Content
$found = $range.Find.Execute("Datasets. (n.d.). Retrieved March 28, 2017, from http://www.civicapps.org/datasets ")
if (-not $found) { throw "Could not find Datasets. anchor text" }
$range.Collapse(1)  # wdCollapseStart

$blob1 = "Databases, Tables & Calculators by Subject. (n.d.). Retrieved April 13, 2017, from https://data.bls.gov/pdq/SurveyOutputServlet " + $cr + $cr + "Bureau of Labor Statistics source for unemployment information. " + $cr + "Series Id: LAUMT413890000000003,LAUMT413890000000004,LAUMT413890000000005,LAUMT413890000000006 " + $cr + "Not Seasonally Adjusted " + $cr + "Area: Portland-Vancouver-Hillsboro, OR-WA Metropolitan Statistical Area" + $cr + "Area Type:Metropolitan areas" + $cr + "State/Region/Division: Oregon" + $cr
$range.InsertBefore($blob1)

# The "Datasets." text now starts a fresh paragraph; find it again so the
# following style/index math for the newly-created paragraphs is reliable.
$range2 = $d.Content
$range2.Find.Execute("Datasets. (n.d.). Retrieved March 28, 2017, from http://www.civicapps.org/datasets ") | Out-Null
$datasetsParaIndex = $range2.Paragraphs(1).Index

# The six "Body Text" paragraphs sit right before the "Datasets." paragraph.
for ($i = $datasetsParaIndex - 6; $i -le $datasetsParaIndex - 1; $i++) {
    $p = $d.Paragraphs($i)
    $p.Style = "Body Text"
    $p.Range.ParagraphFormat.LeftIndent = 35.45
    $p.Range.ParagraphFormat.FirstLineIndent = -0.001
}

# ---------------------------------------------------------------------
# Edit 2: insert the Portland population-history citation + annotation
# right before the "Portland State Criminal Justice Policy Research
# Institute" paragraph.
# ---------------------------------------------------------------------
$range3 = $d.Content
$found3 = $range3.Find.Execute("Portland State Criminal Justice Policy Research Institute")
if (-not $found3) { throw "Could not find Portland State anchor text" }
$range3.Collapse(1)  # wdCollapseStart

$dash = [char]0x2013
$blob2 = "Portland, Oregon Population History 1890 - 2015. (n.d.). Retrieved April 13, 2017, from https://www.biggestuscities.com/city/portland-oregon" + $cr + $cr + "Provides a list of population levels from 1890 " + $dash + " 2015. Decade-level information is available until 2000, then information by year is recorded. Of interest are the 2004-2014 numbers." + $cr + $cr
$range3.InsertBefore($blob2)

$range4 = $d.Content
$range4.Find.Execute("Portland State Criminal Justice Policy Research Institute") | Out-Null
$portlandStateIndex = $range4.Paragraphs(1).Index

# "Provides a list of population levels..." is the 2nd paragraph before
# "Portland State ..." and needs the 709-twip left indent.
$pProvides = $d.Paragraphs($portlandStateIndex - 2)
$pProvides.Range.ParagraphFormat.LeftIndent = 35.45
$pProvides.Range.ParagraphFormat.FirstLineIndent = -0.001

# ---------------------------------------------------------------------
# Edit 3: append the BLS recession-article citation + annotation at the
# very end of the document.
# ---------------------------------------------------------------------
$lastParaIndex = $d.Paragraphs.Count
$lastRange = $d.Paragraphs($lastParaIndex).Range
$lastRange.Collapse(0)  # wdCollapseEnd

$blob3 = $cr + $cr + "U" + "S Bureau of Labor Statistics" + " (2012, February). The Recession of 2007 - 2009. Retrieved April 13, 2017, from https://www.bls.gov/spotlight/2012/recession/pdf/recession_bls_spotlight.pdf " + $cr + $cr + "Article published by the BLS regarding the Recession of 2007 " + $dash + " 2009. Interesting correlations between unemployment and crime are a possibility."
$lastRange.InsertAfter($blob3)

$newLastIndex = $d.Paragraphs.Count
$pArticle = $d.Paragraphs($newLastIndex)
$pArticle.Range.ParagraphFormat.LeftIndent = 35.45
$pArticle.Range.ParagraphFormat.FirstLineIndent = -0.001

Write-Host "Edit complete. Paragraph count:" $d.Paragraphs.Count
